$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.853.52'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.90%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.794.84'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.56%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9996'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.17%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5366'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.41%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3841'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.51%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07444'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.50'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.087'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.71%  '
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.211'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.445'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.03%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.36'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.25%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.789.46'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.99%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.46'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001059'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06523'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.06%  '
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.33'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.966'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.862.88'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.94%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.13'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.092'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.65%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.67'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.32%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.21'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.001.03'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.69%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.331'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '121.54'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.109'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1095'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.99%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.648'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.518'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.06993'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.45%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2190'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.14%  '
$ws.Range("E37").Value = '  -1.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.053'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '11.39'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.84%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.421'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.30%  '
$ws.Range("E41").Value = '  -2.46%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.163'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.22%  '
$ws.Range("B43").Value = 'WEMIXTOKEN'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.410'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.00%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.26'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.678'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5704'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '125.40'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.909'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.172'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.76%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06790'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.32'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.82%  '
